$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2099.5715
$ws.Cells.Item(17, 10).Value = 2099.5715
$ws.Cells.Item(17, 12).Value = 6298.7145
$ws.Cells.Item(17, 14).Value = -6634.7145

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 446.36365
$ws.Cells.Item(53, 9).Value = 345.33334
$ws.Cells.Item(53, 10).Value = 901
$ws.Cells.Item(53, 11).Value = 345.33334
$ws.Cells.Item(53, 12).Value = 901
$ws.Cells.Item(53, 13).Value = 291.66666
$ws.Cells.Item(53, 14).Value = -2175

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2000
$ws.Cells.Item(70, 9).Value = 2000
$ws.Cells.Item(70, 11).Value = 6000
$ws.Cells.Item(70, 13).Value = -5730

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2000
$ws.Cells.Item(73, 9).Value = 2000
$ws.Cells.Item(73, 11).Value = 6000
$ws.Cells.Item(73, 13).Value = -5064

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 4215.1665
$ws.Cells.Item(125, 9).Value = 3869.7144
$ws.Cells.Item(125, 11).Value = 34827.4296
$ws.Cells.Item(125, 13).Value = -32367.4296

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2538.7058
$ws.Cells.Item(137, 9).Value = 1413.7084
$ws.Cells.Item(137, 11).Value = 4241.1252
$ws.Cells.Item(137, 13).Value = -1691.1252

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 8000
$ws.Cells.Item(141, 9).Value = 8000
$ws.Cells.Item(141, 11).Value = 24000
$ws.Cells.Item(141, 13).Value = -18820

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4295.645
$ws.Cells.Item(32, 9).Value = 3938.8333
$ws.Cells.Item(32, 11).Value = 3938.8333
$ws.Cells.Item(32, 13).Value = -3651.8333

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 58040
$ws.Cells.Item(44, 10).Value = 58040
$ws.Cells.Item(44, 12).Value = 58040
$ws.Cells.Item(44, 14).Value = -59016

# ARM row 51
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 18177590
$ws.Cells.Item(74, 9).Value = 33321416
$ws.Cells.Item(74, 10).Value = 4999.6
$ws.Cells.Item(74, 11).Value = 33321416
$ws.Cells.Item(74, 12).Value = 4999.6
$ws.Cells.Item(74, 13).Value = -33320542
$ws.Cells.Item(74, 14).Value = -6747.6

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 18177590
$ws.Cells.Item(77, 9).Value = 33321416
$ws.Cells.Item(77, 10).Value = 4999.6
$ws.Cells.Item(77, 11).Value = 166607080
$ws.Cells.Item(77, 12).Value = 24998
$ws.Cells.Item(77, 13).Value = -166602712
$ws.Cells.Item(77, 14).Value = -33734

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1331.7693
$ws.Cells.Item(122, 9).Value = 1331.7693
$ws.Cells.Item(122, 11).Value = 3995.3079
$ws.Cells.Item(122, 13).Value = -1545.3079

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4335.5625
$ws.Cells.Item(86, 9).Value = 3402.4443
$ws.Cells.Item(86, 10).Value = 5535.2856
$ws.Cells.Item(86, 11).Value = 3402.4443
$ws.Cells.Item(86, 12).Value = 5535.2856
$ws.Cells.Item(86, 13).Value = -2279.4443
$ws.Cells.Item(86, 14).Value = -7781.2856

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 4335.5625
$ws.Cells.Item(89, 9).Value = 3402.4443
$ws.Cells.Item(89, 10).Value = 5535.2856
$ws.Cells.Item(89, 11).Value = 17012.2215
$ws.Cells.Item(89, 12).Value = 27676.428
$ws.Cells.Item(89, 13).Value = -11396.2215
$ws.Cells.Item(89, 14).Value = -38908.428

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1553.5625
$ws.Cells.Item(107, 9).Value = 1234.6364
$ws.Cells.Item(107, 11).Value = 1234.6364
$ws.Cells.Item(107, 13).Value = 685.3635999999999

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 520
$ws.Cells.Item(7, 10).Value = 574.75
$ws.Cells.Item(7, 12).Value = 574.75
$ws.Cells.Item(7, 14).Value = -800.75

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2845.6667
$ws.Cells.Item(31, 10).Value = 2014
$ws.Cells.Item(31, 12).Value = 2014
$ws.Cells.Item(31, 14).Value = -2604

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2845.6667
$ws.Cells.Item(34, 10).Value = 2014
$ws.Cells.Item(34, 12).Value = 2014
$ws.Cells.Item(34, 14).Value = -2418

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1515.9166
$ws.Cells.Item(134, 9).Value = 1562.8182
$ws.Cells.Item(134, 11).Value = 4688.4546
$ws.Cells.Item(134, 13).Value = -2153.4546

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 90000
$ws.Cells.Item(141, 10).Value = 85000
$ws.Cells.Item(141, 12).Value = 85000
$ws.Cells.Item(141, 14).Value = -95360

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 2625.75
$ws.Cells.Item(46, 9).Value = 2201
$ws.Cells.Item(46, 10).Value = 3900
$ws.Cells.Item(46, 11).Value = 6603
$ws.Cells.Item(46, 12).Value = 11700
$ws.Cells.Item(46, 13).Value = -6512
$ws.Cells.Item(46, 14).Value = -11882

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 253293
$ws.Cells.Item(80, 9).Value = 4391.5
$ws.Cells.Item(80, 10).Value = 502194.5
$ws.Cells.Item(80, 11).Value = 13174.5
$ws.Cells.Item(80, 12).Value = 1506583.5
$ws.Cells.Item(80, 13).Value = -12238.5
$ws.Cells.Item(80, 14).Value = -1508455.5

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 253293
$ws.Cells.Item(83, 9).Value = 4391.5
$ws.Cells.Item(83, 10).Value = 502194.5
$ws.Cells.Item(83, 11).Value = 39523.5
$ws.Cells.Item(83, 12).Value = 4519750.5
$ws.Cells.Item(83, 13).Value = -34843.5
$ws.Cells.Item(83, 14).Value = -4529110.5

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 1579.6
$ws.Cells.Item(87, 9).Value = 1579.6
$ws.Cells.Item(87, 11).Value = 4738.799999999999
$ws.Cells.Item(87, 13).Value = -3490.799999999999

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 1579.6
$ws.Cells.Item(90, 9).Value = 1579.6
$ws.Cells.Item(90, 11).Value = 14216.4
$ws.Cells.Item(90, 13).Value = -7976.4

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 587.25
$ws.Cells.Item(92, 9).Value = 299
$ws.Cells.Item(92, 11).Value = 897
$ws.Cells.Item(92, 13).Value = 351

# CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 625.6667
$ws.Cells.Item(103, 9).Value = 724.5
$ws.Cells.Item(103, 10).Value = 428
$ws.Cells.Item(103, 11).Value = 2173.5
$ws.Cells.Item(103, 12).Value = 1284
$ws.Cells.Item(103, 13).Value = -1294.5
$ws.Cells.Item(103, 14).Value = -3042

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1371
$ws.Cells.Item(131, 10).Value = 2331
$ws.Cells.Item(131, 12).Value = 6993
$ws.Cells.Item(131, 14).Value = -17073

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 5254.8
$ws.Cells.Item(138, 9).Value = 4318.5
$ws.Cells.Item(138, 11).Value = 12955.5
$ws.Cells.Item(138, 13).Value = -7815.5

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1381.1538
$ws.Cells.Item(122, 9).Value = 1246
$ws.Cells.Item(122, 10).Value = 1831.6666
$ws.Cells.Item(122, 11).Value = 3738
$ws.Cells.Item(122, 12).Value = 5494.9998
$ws.Cells.Item(122, 13).Value = -1288
$ws.Cells.Item(122, 14).Value = -10394.9998

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2895
$ws.Cells.Item(132, 9).Value = 1953.6666
$ws.Cells.Item(132, 10).Value = 3459.8
$ws.Cells.Item(132, 11).Value = 5860.9998
$ws.Cells.Item(132, 12).Value = 10379.4
$ws.Cells.Item(132, 13).Value = -3330.9998
$ws.Cells.Item(132, 14).Value = -15439.4

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 939.2
$ws.Cells.Item(46, 9).Value = 924.25
$ws.Cells.Item(46, 11).Value = 924.25
$ws.Cells.Item(46, 13).Value = -736.25

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 298.75
$ws.Cells.Item(55, 9).Value = 298.5
$ws.Cells.Item(55, 10).Value = 299.5
$ws.Cells.Item(55, 11).Value = 298.5
$ws.Cells.Item(55, 12).Value = 299.5
$ws.Cells.Item(55, 13).Value = -125.5
$ws.Cells.Item(55, 14).Value = -645.5

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1475.75
$ws.Cells.Item(82, 9).Value = 1200
$ws.Cells.Item(82, 10).Value = 1751.5
$ws.Cells.Item(82, 11).Value = 1200
$ws.Cells.Item(82, 12).Value = 1751.5
$ws.Cells.Item(82, 13).Value = -839
$ws.Cells.Item(82, 14).Value = -2473.5

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1475.75
$ws.Cells.Item(85, 9).Value = 1200
$ws.Cells.Item(85, 10).Value = 1751.5
$ws.Cells.Item(85, 11).Value = 1200
$ws.Cells.Item(85, 12).Value = 1751.5
$ws.Cells.Item(85, 13).Value = 48
$ws.Cells.Item(85, 14).Value = -4247.5

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5332.3335
$ws.Cells.Item(62, 9).Value = 5332.3335
$ws.Cells.Item(62, 11).Value = 5332.3335
$ws.Cells.Item(62, 13).Value = -4708.3335

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 5332.3335
$ws.Cells.Item(65, 9).Value = 5332.3335
$ws.Cells.Item(65, 11).Value = 26661.6675
$ws.Cells.Item(65, 13).Value = -23541.6675
